$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.058859333333333
$ws.Range("H2").Value = 6.176577999999999
$ws.Range("I2").Value = 0.03050820259458848
$ws.Range("J2").Value = 0.03050820259458848
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 158.7164846170511
$ws.Range("R2").Value = 1428.44836155346
$ws.Range("S2").Value = 0.007333612608230775
$ws.Range("T2").Value = 0.007333612608230775
$ws.Range("G3").Value = 2.058859333333333
$ws.Range("H3").Value = 6.176577999999999
$ws.Range("I3").Value = 0.03050820259458848
$ws.Range("J3").Value = 0.03050820259458848
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 209.1390079440818
$ws.Range("R3").Value = 1882.251071496736
$ws.Range("S3").Value = 0.009663422606871563
$ws.Range("T3").Value = 0.009663422606871561
$ws.Range("G4").Value = 2.058859333333333
$ws.Range("H4").Value = 6.176577999999999
$ws.Range("I4").Value = 0.03050820259458848
$ws.Range("J4").Value = 0.03050820259458848
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 292.4131808023004
$ws.Range("R4").Value = 2631.718627220704
$ws.Range("S4").Value = 0.01351116737948614
$ws.Range("T4").Value = 0.01351116737948614
$ws.Range("I5").Value = 0.540047065760451
$ws.Range("J5").Value = 0.540047065760451
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 2809.551678421597
$ws.Range("R5").Value = 25285.96510579438
$ws.Range("S5").Value = 0.1298174141272219
$ws.Range("T5").Value = 0.1298174141272219
$ws.Range("I6").Value = 0.540047065760451
$ws.Range("J6").Value = 0.540047065760451
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.1710590129937673
$ws.Range("T6").Value = 0.1710590129937673
$ws.Range("I7").Value = 0.540047065760451
$ws.Range("J7").Value = 0.540047065760451
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.2391706386394619
$ws.Range("T7").Value = 0.2391706386394619
$ws.Range("I8").Value = 0.4294447316449605
$ws.Range("J8").Value = 0.4294447316449605
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 2234.151878750509
$ws.Range("R8").Value = 20107.36690875458
$ws.Range("S8").Value = 0.1032306406372298
$ws.Range("T8").Value = 0.1032306406372298
$ws.Range("I9").Value = 0.4294447316449605
$ws.Range("J9").Value = 0.4294447316449605
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.1360259069774209
$ws.Range("T9").Value = 0.1360259069774208
$ws.Range("I10").Value = 0.4294447316449605
$ws.Range("J10").Value = 0.4294447316449605
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.1901881840303099
$ws.Range("T10").Value = 0.1901881840303099
